# Regenerate merged AHB files
# 1. Rename the header row: "<Feld>_old" -> "<Feld>_FV2210" and "<Feld>_new" -> "<Feld>_FV2304"
# 2. Turn the data range into an Excel Table ("Table1")
# 3. Freeze the header row (split/freeze at row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2210",
  "Segmentgruppe_FV2210",
  "Segment_FV2210",
  "Datenelement_FV2210",
  "Segment ID_FV2210",
  "Code_FV2210",
  "Qualifier_FV2210",
  "Beschreibung_FV2210",
  "Bedingungsausdruck_FV2210",
  "Bedingung_FV2210",
  "diff",
  "Segmentname_FV2304",
  "Segmentgruppe_FV2304",
  "Segment_FV2304",
  "Datenelement_FV2304",
  "Segment ID_FV2304",
  "Code_FV2304",
  "Qualifier_FV2304",
  "Beschreibung_FV2304",
  "Bedingungsausdruck_FV2304",
  "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn A1:U94 into a real Excel table
$rng = $ws.Range("A1:U94")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
